# The SMOTE / cluster-analysis re-run produced columns E..I in the wrong
# order. Column E should hold what is currently in G, F should hold what
# is currently in H, G should hold what is currently in E, H should hold
# what is currently in I, and I should hold what is currently in F:
#
#   new E = old G
#   new F = old H
#   new G = old E
#   new H = old I
#   new I = old F
#
# That permutation decomposes into two independent cycles:
#   (E G)       -- a plain swap
#   (F I H)     -- F -> I -> H -> F
#
# We rotate each cycle with Copy / PasteSpecial(xlPasteValues) so that only
# the cell VALUE moves — formatting/style of every destination cell is left
# completely untouched, matching how the column mix-up was actually fixed.
# A1 (always blank in this sheet) is reused as the one-cell scratch buffer
# needed to rotate a cycle without clobbering a value before it's been
# copied onward; it is restored to blank once every row has been fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteValues = -4163
$scratch = "A1"

$lastRow = 71   # sheet's populated range is A1:I71

for ($r = 1; $r -le $lastRow; $r++) {

    # --- cycle (E G): swap E<r> and G<r> ------------------------------
    $ws.Range("G$r").Copy()
    $ws.Range($scratch).PasteSpecial($xlPasteValues)
    $ws.Range("E$r").Copy()
    $ws.Range("G$r").PasteSpecial($xlPasteValues)
    $ws.Range($scratch).Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteValues)

    # --- cycle (F I H): rotate F<r> -> I<r> -> H<r> -> F<r> -----------
    $ws.Range("H$r").Copy()
    $ws.Range($scratch).PasteSpecial($xlPasteValues)
    $ws.Range("I$r").Copy()
    $ws.Range("H$r").PasteSpecial($xlPasteValues)
    $ws.Range("F$r").Copy()
    $ws.Range("I$r").PasteSpecial($xlPasteValues)
    $ws.Range($scratch).Copy()
    $ws.Range("F$r").PasteSpecial($xlPasteValues)
}

# Leave the scratch cell the way we found it: empty.
$ws.Range($scratch).ClearContents()
